$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new time log entry in row 111
# (Interruption is set before Start/Stop time so the Delta formula's
#  dependency on D111 recalculates correctly alongside B111/C111.)
$ws.Range("A111").Value = 41952
$ws.Range("D111").Value = 15
$ws.Range("B111").Value = 0.52638888888888891
$ws.Range("C111").Value = 0.63680555555555551
$ws.Range("F111").Value = "Coding"

# Update selection to match the final state (B112)
$ws.Range("B112").Select()
